$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently ends with:
#   row 20 : ... Find Subsequence of Length K With the Largest Sum ...
#   row 21 : blank placeholder row (only H21/I21 styled, no values)
#   row 22 : 594  Longest Harmonious Subsequence ...
#
# The edit inserts a brand-new LeetCode #1498 entry where the old blank
# row 21 used to be (pushing the "594" row down to row 22, which it
# already effectively is), and appends a brand-new LeetCode #4 entry as
# the new row 23.
# ---------------------------------------------------------------------------

# Insert a fresh row at 21 (shifts old row21/row22 down by one).
$ws.Rows.Item(21).Insert()

# New row 21: LeetCode 1498 - Number of Subsequences That Satisfy the Given Sum Condition
$ws.Range("A21").Value = 1498
$ws.Range("B21").Value = "Number of Subsequences That Satisfy the Given Sum Condition"
$ws.Range("C21").Value = "#array #two-pointers #binary-search #sorting "
$ws.Range("D21").Value = "medium"
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 18
$ws.Range("H21").Value = 45837
$ws.Range("I21").Value = 45837
$ws.Rows.Item(21).RowHeight = 51

# The old placeholder row (now pushed to row 22, still empty) is no longer
# needed now that row 21 carries real data - drop it so the "594" entry
# lands back on row 22.
$ws.Rows.Item(22).Delete()

# New row 23: LeetCode 4 - Median of Two Sorted Arrays
$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "Median of Two Sorted Arrays"
$ws.Range("C23").Value = "#array #binary-search #divide-and-conquer #核心 "
$ws.Range("D23").Value = "hard"
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 18
$ws.Range("H23").Value = 45838
$ws.Range("I23").Value = 45838
$ws.Rows.Item(23).RowHeight = 51

# Match the date-column formatting (style) already used by row 22's H/I cells.
$ws.Range("H22:I22").Copy()
$ws.Range("H23:I23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the saved selection to match the new cursor position.
$ws.Range("F24").Select()
